$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 783
$ws1.Range("F7").Value = 649
$ws1.Range("F12").Value = 528
$ws1.Range("F13").Value = 174
$ws1.Range("F15").Value = 718
$ws1.Range("F22").Value = 607
$ws1.Range("F24").Value = 869

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 332
$ws2.Range("F4").Value = 107
$ws2.Range("F9").Value = 53

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 332
$ws4.Range("F5").Value = 783
$ws4.Range("F10").Value = 649
$ws4.Range("F14").Value = 107
$ws4.Range("F17").Value = 528
$ws4.Range("F19").Value = 174
$ws4.Range("F21").Value = 718
$ws4.Range("F28").Value = 53
$ws4.Range("F35").Value = 607
$ws4.Range("F37").Value = 869
